$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts the existing B:J data to C:K
$ws.Columns("B:B").Insert()

# Header for the new column
$ws.Range("B1").Value = "_requirements"

# Values for the new column on the rows that need them
$ws.Range("B5").Value = "l10n_it_reverse_charge"
$ws.Range("B6").Value = "l10n_it_split_payment"
$ws.Range("B7").Value = "l10n_it_dichiarazione_intento or l10n_it_lettera_intento"

# Set the new column's width (~44.6 characters, matching the source file)
$ws.Columns("B:B").ColumnWidth = 43.83

# Update selection to match target state
$ws.Range("B8").Select() | Out-Null
